$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("windfarms")

# Add a new data row (row 2) describing the "SLIDER_WF" wind farm entry.
$ws.Range("A2").Value = "SLIDER_WF"
$ws.Range("B2").Value = "SLIDER_WF"
$ws.Range("C2").Value = "SLIDER_WF"
$ws.Range("K2").Value = $true
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 5000
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 5000
$ws.Range("P2").Value = 52
$ws.Range("Q2").Value = 5

# Match the author's final selection as recorded in the saved workbook.
$ws.Range("N12").Select()
